$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Capture the quote-prefix cell style (currently on C30, the row for
# "LE_MF_PICOS_MEMBER_PROXY_APP_B") into a scratch cell BEFORE the row
# insert below, since EntireRow.Insert() does not preserve the
# quote-prefix flag when it shifts a row's formatting down.
# ------------------------------------------------------------------
$ws.Range("C30").Copy()
$ws.Range("ZZ1").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Insert a new row at 25 for the new "LE_WHITELISTED_P_FUND_B" entry;
# this shifts existing rows 25-36 down to 26-37 and all the formulas
# that reference B-column rows (e.g. "=B30+A13") auto-adjust.
# ------------------------------------------------------------------
$ws.Range("A25").EntireRow.Insert()

# New row 25: LE_WHITELISTED_B | LE_P_FUND_B = 32 + 16 = 48
$ws.Range("B25").Formula = "=A7+A6"

# ------------------------------------------------------------------
# Update every label in column C (rows 2-37) to the new text per the
# renamed / renumbered / reworded bit-flag constants.
# ------------------------------------------------------------------
$ws.Range("C2").Value = ' LE_REGISTERED_B                =       1; //  0 Entry has been registered with addedT set but nothing more'
$ws.Range("C3").Value = ' LE_SALE_CONTRACT_B             =       2; //  1 Is the Sale Contract entry - where the minted PIOs are held. Has dbId == 1'
$ws.Range("C4").Value = ' LE_M_FUND_B                    =       4; //  2 Mfund funded whitelisted with picos entry or unfunded whitelisted with picos entry. See below for more.'
$ws.Range("C5").Value = ' LE_PICOS_B                     =       8; //  3 Holds Picos. Can be set wo LE_M_FUND_B being set for a presale entry'
$ws.Range("C6").Value = ' LE_P_FUND_B                    =      16; //  4 Pfund prepurchase entry, always funded. There are 4 types of prepurchase entries as below. If unset then entry is an escrow entry, and must then have either LE_WHITELISTED_B or LE_PRESALE_B set or both.'
$ws.Range("C7").Value = ' LE_WHITELISTED_B               =      32; //  5 Has been whitelisted'
$ws.Range("C8").Value = ' LE_MEMBER_B                    =      64; //  6 Is a Pacio Member: Whitelisted with a picosBalance'
$ws.Range("C9").Value = ' LE_PRESALE_B                   =     128; //  7 A Presale List entry - Pacio Seed Presale or Pacio Private Placement. /- Can make Tranche 1 purchases but not Tranche 2 to 4 ones on same account'
$ws.Range("C10").Value = ' LE_TRANCH1_B                   =     256; //  8 Was or included a Tranche 1 purchase.                                 |   until after soft cap as not entitled to soft cap miss refund                              -'
$ws.Range("C11").Value = ' LE_FROM_TRANSFER_OK_B          =     512; //  9 Transfers from this entry allowed entry even if pTransfersOkB is false. Is set for the Sale contract entry.'
$ws.Range("C12").Value = ' LE_PROXY_APPOINTER_B           =    1024; // 10 This entry has appointed a Proxy. Need not be a Member.                                              /- one entry can have both bits set'
$ws.Range("C13").Value = ' LE_PROXY_B                     =    2048; // 11 This entry is a Proxy i.e. one or more other entries have appointed it as a proxy. Must be a Member. |  as a proxy can appoint a proxy'
$ws.Range("C14").Value = ' LE_DOWNGRADED_B                =    4096; // 12 This entry has been downgraded from whitelisted. Refunding candidate.'
$ws.Range("C15").Value = ' LE_BLOCKED_FROM_VOTING_B       =    8192; // 13 Set if a member is blocked from voting by a PGC managed op as a result of trolling etc'
$ws.Range("C16").Value = ' LE_TRANSFERRED_TO_PB_B         =   16384; // 14 This entry has had its PIOs transferred to the Pacio Blockchain'
$ws.Range("C17").Value = ' LE_P_REFUNDED_S_CAP_MISS_B     =   32768; // 15 Pfund funds Refunded due to soft cap not being reached'
$ws.Range("C18").Value = ' LE_P_REFUNDED_SALE_CLOSE_B     =   65536; // 16 Pfund funds Refunded due to not being whitelisted by the time that the sale closes'
$ws.Range("C19").Value = ' LE_P_REFUNDED_ONCE_OFF_B       =  131072; // 17 Pfund funds Refunded once off manually for whatever reason'
$ws.Range("C20").Value = ' LE_M_REFUNDED_S_CAP_MISS_NPT1B =  262144; // 18 Mfund funds Refunded due to soft cap not being reached. Such refunds do not apply to MFunds from a presale or tranche 1 purchase.'
$ws.Range("C21").Value = ' LE_M_REFUNDED_TERMINATION_B    =  524288; // 19 Mfund or Presale with picos Refund proportionately according to Picos held following a vote for project termination'
$ws.Range("C22").Value = ' LE_M_REFUNDED_ONCE_OFF_B       = 1048576; // 20 Mfund funds Refunded once off manually for whatever reason including downgrade from whitelisted'
$ws.Range("C23").Value = ' // Combos'
$ws.Range("C24").Value = ' LE_M_FUND_PICOS_MEMBER_B       =   76; // LE_M_FUND_B | LE_PICOS_B | LE_MEMBER_B'
$ws.Range("C25").Value = ' LE_WHITELISTED_P_FUND_B        =   48; // LE_WHITELISTED_B | LE_P_FUND_B'
$ws.Range("C26").Value = ' LE_WHITELISTED_MEMBER_B        =   96; // LE_WHITELISTED_B | LE_MEMBER_B'
$ws.Range("C27").Value = ' LE_PRESALE_TRANCH1_B           =  384; // LE_PRESALE_B | LE_TRANCH1_B == not eligible for a soft cap miss refund'
$ws.Range("C28").Value = ' LE_MEMBER_PROXY_B              = 2112; // LE_MEMBER_B | LE_PROXY_B'
$ws.Range("C29").Value = ' LE_PROXY_INVOLVED_COMBO_B      = 3072; // LE_PROXY_APPOINTER_B | LE_PROXY_B'
$ws.Range("C30").Value = ' LE_PROXY_APP_VOTE_BLOCK_B      = 9216; // LE_PROXY_APPOINTER_B | LE_BLOCKED_FROM_VOTING_B'
$ws.Range("C31").Value = ' LE_MF_PICOS_MEMBER_PROXY_APP_B = 1100; // LE_M_FUND_B | LE_PICOS_B | LE_MEMBER_B | LE_PROXY_APPOINTER_B'
$ws.Range("C32").Value = ' LE_MF_PICOS_MEMBER_PROXY_ALL_B = 3148; // LE_M_FUND_B | LE_PICOS_B | LE_MEMBER_B | LE_PROXY_INVOLVED_COMBO_B'
$ws.Range("C33").Value = ' LE_REFUNDED_COMBO_B         = 2064384; // LE_P_REFUNDED_S_CAP_MISS_B | LE_P_REFUNDED_SALE_CLOSE_B | LE_P_REFUNDED_ONCE_OFF_B | LE_M_REFUNDED_S_CAP_MISS_NPT1B | LE_M_REFUNDED_TERMINATION_B | LE_M_REFUNDED_ONCE_OFF_B'
$ws.Range("C34").Value = ' LE_DEAD_COMBO_B             = 2080768; // LE_TRANSFERRED_TO_PB_B | LE_REFUNDED_COMBO_B  or bits >= 8192'
$ws.Range("C35").Value = ' LE_NO_SEND_FUNDS_COMBO_B    = 2084994; // LE_DEAD_COMBO_B | LE_SALE_CONTRACT_B | LE_PRESALE | LE_DOWNGRADED_B'
$ws.Range("C36").Value = ' LE_NO_REFUNDS_COMBO_B       = 2080770; // LE_DEAD_COMBO_B | LE_SALE_CONTRACT_B Starting point check. Could also be more i.e. no funds or no PIOs'

# ------------------------------------------------------------------
# Restore the quote-prefix style onto C31 (the row that used to be
# C30 pre-insert) and clean up the scratch cell.
# ------------------------------------------------------------------
$ws.Range("ZZ1").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("ZZ1").Clear()

# Match the saved selection state recorded in the workbook.
$ws.Range("D34").Select() | Out-Null
